# Actualizacion automatica del tracker
# Append the latest result row (row 60) to the bottom of the tracker table
# on Sheet1, following the same A:H layout as every existing row
# (event_id, fecha, jugador_A, jugador_B, pronostico, cuota, resultado,
# profit). This particular match hasn't been settled yet, so "resultado"
# and "profit" (columns G/H) are left blank, same as the other still-open
# rows further up the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 60

$ws.Cells.Item($row, 1).Value = 14739135
$ws.Cells.Item($row, 2).Value = "'2025-09-22"
$ws.Cells.Item($row, 2).Style = "Normal"
$ws.Cells.Item($row, 3).Value = "Gonzalo Villanueva"
$ws.Cells.Item($row, 4).Value = "Facundo Juarez"
$ws.Cells.Item($row, 5).Value = "Gana Facundo Juarez"
$ws.Cells.Item($row, 6).Value = 2.5
$ws.Cells.Item($row, 7).Value = "'"
$ws.Cells.Item($row, 7).Style = "Normal"
$ws.Cells.Item($row, 8).Value = "'"
$ws.Cells.Item($row, 8).Style = "Normal"
